$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: RowNumber, Coin(B), Link(C), Price(D), Volume1h(E)
$rows = @(
    @(2, 'Bitcoin', 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc', '57.364.17', '  +2.43%  '),
    @(3, 'Ethereum', 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth', '2.346.89', '  +1.51%  '),
    @(4, 'TetherUSD', 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt', '1.00', '  -0.13%  '),
    @(5, 'BNB', 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb', '520.94', '  +0.77%  '),
    @(6, 'Solana', 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol', '136.95', '  +3.15%  '),
    @(7, 'USDC', 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc', '0.997', '  +0.23%  '),
    @(8, 'XRP', 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp', '0.538', '  +0.97%  '),
    @(9, 'LidoStakedEther', 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth', '2.357.95', '  +1.14%  '),
    @(10, 'Dogecoin', 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge', '0.103', '  -0.36%  '),
    @(11, 'Toncoin', 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton', '5.47', '  +6.38%  '),
    @(12, 'TRON', 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx', '0.152', '  -1.26%  '),
    @(13, 'Cardano', 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada', '0.344', '  +0.29%  '),
    @(14, 'Avalanche', 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax', '24.13', '  +0.50%  '),
    @(15, 'WrappedliquidstakedEther2.0', 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth', '2.760.64', '  +1.42%  '),
    @(16, 'WrappedBTC', 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc', '57.316.99', '  +2.03%  '),
    @(17, 'ShibaInu', 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib', '0.0000136', '  +0.23%  '),
    @(18, 'WrappedEther', 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth', '2.346.70', '  +0.20%  '),
    @(19, 'Chainlink', 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link', '10.66', '  +1.61%  '),
    @(20, 'BitcoinCash', 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch', '329.06', '  +2.32%  '),
    @(21, 'Polkadot', 'https://coinranking.com/coin/25W7FG7om+polkadot-dot', '4.25', '  -0.64%  '),
    @(22, 'Uniswap', 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni', '6.75', '  +1.34%  '),
    @(23, 'Dai', 'https://coinranking.com/coin/MoTuySvg7+dai-dai', '0.999', '  +0.01%  '),
    @(24, 'Litecoin', 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc', '61.58', '  +1.54%  '),
    @(25, 'Kaspa', 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas', '0.166', '  +5.08%  '),
    @(26, 'Binance-PegBSC-USD', 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd', '0.995', '  +0.40%  '),
    @(27, 'InternetComputer(DFINITY)', 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp', '8.23', '  +7.51%  '),
    @(28, 'Fetch.AI', 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet', '1.32', '  +10.96%  '),
    @(29, 'PEPE', 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe', '0.0₃0746', '  +2.69%  '),
    @(30, 'Monero', 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr', '170.07', '  -1.10%  '),
    @(31, 'PancakeSwap', 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake', '1.71', '  +1.55%  '),
    @(32, 'Aptos', 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt', '6.26', '  +0.15%  '),
    @(33, 'EthereumClassic', 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc', '18.64', '  +1.68%  '),
    @(34, 'USDe', 'https://coinranking.com/coin/exbfr2U-0+usde-usde', '0.999', '  +0.04%  '),
    @(35, 'ImmutableX', 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx', '1.28', '  +1.74%  '),
    @(36, 'FirstDigitalUSD', 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd', '0.994', '  +0.13%  '),
    @(37, 'SuiNetwork', 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui', '0.924', '  -0.32%  '),
    @(38, 'NEARProtocol', 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near', '4.05', '  +1.59%  '),
    @(39, 'Stacks', 'https://coinranking.com/coin/mMPrMcB7+stacks-stx', '1.58', '  +4.58%  '),
    @(40, 'OKB', 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb', '38.54', '  +3.19%  '),
    @(41, 'Aave', 'https://coinranking.com/coin/ixgUfzmLR+aave-aave', '149.36', '  +7.14%  '),
    @(42, 'PolygonEcosystemToken', 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol', '0.385', '  +0.32%  '),
    @(43, 'Filecoin', 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil', '3.65', '  +2.30%  '),
    @(44, 'Bittensor', 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao', '283.76', '  +5.03%  '),
    @(45, 'RenderToken', 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr', '5.23', '  +3.56%  '),
    @(46, 'Stellar', 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm', '0.0938', '  +1.34%  '),
    @(47, 'Hedera', 'https://coinranking.com/coin/jad286TjB+hedera-hbar', '0.0507', '  -0.46%  '),
    @(48, 'Mantle', 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt', '0.564', '  +2.14%  '),
    @(49, 'VeChain', 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet', '0.0221', '  +2.63%  '),
    @(50, 'InjectiveProtocol', 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj', '18.21', '  +6.51%  '),
    @(51, 'EnergySwap', 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens', '17.57', '  +4.03%  ')
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).NumberFormat = "@"
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).NumberFormat = "@"
    $ws.Cells.Item($r, 5).Value = $row[4]
}
